# 2014 Metadata links updated
#
# Applies the "2014 Metadata links updated" commit to the RP1 SAF RAT
# application workbook:
#  - The "Meta data" link (I1) on the "RAT 2014" and "RAT EU_2014" sheets
#    no longer shows the raw wiki page title; it now reads
#    "Metadata - Single European Sky Portal" (keeping the same hyperlink
#    target).
#  - The "Contact" link (I2) on both sheets no longer reads the NSA
#    address as a mailto hyperlink; it is now the plain support address
#    "pru-support@eurocontrol.int" (hyperlink removed).
#  - A data correction: the Malta "ATM tech. events (%ATM Overall)" figure
#    on "RAT 2014" (H23) is corrected from 0.06 to 1.0, and the EU-wide
#    2014 "ATM tech. events (overall)" figure on "RAT EU_2014" (H8) is
#    updated from 0.7 to 0.73.
#  - A new "Change Log" entry is added documenting the Malta correction.

$wb = $excel.ActiveWorkbook

$newMetaText = "Metadata - Single European Sky Portal"
$newContactText = "pru-support@eurocontrol.int"

foreach ($sheetName in @("RAT 2014", "RAT EU_2014")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove only the mailto: hyperlink (on I2); keep the wiki link on I1.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Address -like "mailto:*") {
            $h.Delete()
        }
    }

    # Replace the formula-driven HYPERLINK() text with plain values,
    # preserving each cell's existing style.
    $ws.Range("I1").Value = $newMetaText
    $ws.Range("I2").Value = $newContactText
}

# Data correction: Malta ATM tech. events (%ATM Overall) on "RAT 2014".
$wsRat2014 = $wb.Worksheets.Item("RAT 2014")
$wsRat2014.Range("H23").Value = 1.0

# Data correction: EU-wide 2014 ATM tech. events (overall) on "RAT EU_2014".
$wsRatEu2014 = $wb.Worksheets.Item("RAT EU_2014")
$wsRatEu2014.Range("H8").Value = 0.73

# Append the corresponding Change Log entry, matching the formatting of
# the preceding row.
$wsLog = $wb.Worksheets.Item("Change Log")
$wsLog.Range("A29:D29").Copy($wsLog.Range("A30:D30"))
$wsLog.Range("A30").Value = "4/9/2016"
$wsLog.Range("B30").Value = "Malta"
$wsLog.Range("C30").Value = 2014
$wsLog.Range("D30").Value = "ATM specific tech. events updated"
$wsLog.Rows.Item(30).RowHeight = 12.75
